# Applies the 2023-11-05 08:45 scrape update to the Germany Bundesliga 2023-2024 sheet:
#  - re-orders the F:V (match) columns within several same-matchday row blocks
#    (the A-E "index/date" columns already line up; only match/odds data moved)
#  - appends 6 newly scraped matches as rows 84-89

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: capture current F:V values for every row that will be rewritten ---
# (read everything first so the subsequent writes cannot clobber a value we still need)
$src39 = $ws.Range("F39:V39").Value()
$src40 = $ws.Range("F40:V40").Value()
$src41 = $ws.Range("F41:V41").Value()
$src42 = $ws.Range("F42:V42").Value()
$src48 = $ws.Range("F48:V48").Value()
$src49 = $ws.Range("F49:V49").Value()
$src50 = $ws.Range("F50:V50").Value()
$src51 = $ws.Range("F51:V51").Value()
$src52 = $ws.Range("F52:V52").Value()
$src57 = $ws.Range("F57:V57").Value()
$src58 = $ws.Range("F58:V58").Value()
$src59 = $ws.Range("F59:V59").Value()
$src60 = $ws.Range("F60:V60").Value()
$src66 = $ws.Range("F66:V66").Value()
$src67 = $ws.Range("F67:V67").Value()
$src68 = $ws.Range("F68:V68").Value()
$src69 = $ws.Range("F69:V69").Value()
$src70 = $ws.Range("F70:V70").Value()
$src75 = $ws.Range("F75:V75").Value()
$src76 = $ws.Range("F76:V76").Value()
$src77 = $ws.Range("F77:V77").Value()
$src78 = $ws.Range("F78:V78").Value()
$src79 = $ws.Range("F79:V79").Value()

# --- Step 2: write the captured values back out to their new row positions ---
$ws.Range("F39:V39").Value = $src41
$ws.Range("F40:V40").Value = $src39
$ws.Range("F41:V41").Value = $src42
$ws.Range("F42:V42").Value = $src40
$ws.Range("F48:V48").Value = $src51
$ws.Range("F49:V49").Value = $src52
$ws.Range("F50:V50").Value = $src49
$ws.Range("F51:V51").Value = $src50
$ws.Range("F52:V52").Value = $src48
$ws.Range("F57:V57").Value = $src60
$ws.Range("F58:V58").Value = $src59
$ws.Range("F59:V59").Value = $src57
$ws.Range("F60:V60").Value = $src58
$ws.Range("F66:V66").Value = $src70
$ws.Range("F67:V67").Value = $src69
$ws.Range("F68:V68").Value = $src66
$ws.Range("F69:V69").Value = $src67
$ws.Range("F70:V70").Value = $src68
$ws.Range("F75:V75").Value = $src79
$ws.Range("F76:V76").Value = $src78
$ws.Range("F77:V77").Value = $src75
$ws.Range("F78:V78").Value = $src76
$ws.Range("F79:V79").Value = $src77

# --- Step 3: append the 6 newly scraped matches as rows 84-89 ---
# Copy formatting (style/number format/border) from the last existing data row (83)
# across the new block, then fill in the actual values.
$ws.Range("A83:V83").Copy($ws.Range("A84:V89"))

# Row 84
$ws.Range("A84").Value = 83
$ws.Range("B84").Value = "germany"
$ws.Range("C84").Value = "bundesliga"
$ws.Range("D84").Value = "2023-2024"
$ws.Range("E84").Value = 45234.64583333334
$ws.Range("F84").Value = "Mainz"
$ws.Range("G84").Value = 2
$ws.Range("H84").Value = "RB Leipzig"
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 3.62
$ws.Range("K84").Value = "22/10/2023 12:02"
$ws.Range("L84").Value = 4.57
$ws.Range("M84").Value = "04/11/2023 15:28"
$ws.Range("N84").Value = 3.78
$ws.Range("O84").Value = "22/10/2023 12:02"
$ws.Range("P84").Value = 4.04
$ws.Range("Q84").Value = "04/11/2023 15:28"
$ws.Range("R84").Value = 1.93
$ws.Range("S84").Value = "22/10/2023 12:02"
$ws.Range("T84").Value = 1.78
$ws.Range("U84").Value = "04/11/2023 15:28"
$ws.Range("V84").Value = "https://www.betexplorer.com/football/germany/bundesliga/mainz-rb-leipzig/r1DeDG8e/"

# Row 85
$ws.Range("A85").Value = 84
$ws.Range("B85").Value = "germany"
$ws.Range("C85").Value = "bundesliga"
$ws.Range("D85").Value = "2023-2024"
$ws.Range("E85").Value = 45234.64583333334
$ws.Range("F85").Value = "FC Koln"
$ws.Range("G85").Value = 1
$ws.Range("H85").Value = "Augsburg"
$ws.Range("I85").Value = 1
$ws.Range("J85").Value = 1.9
$ws.Range("K85").Value = "22/10/2023 12:02"
$ws.Range("L85").Value = 2.05
$ws.Range("M85").Value = "04/11/2023 15:29"
$ws.Range("N85").Value = 3.73
$ws.Range("O85").Value = "22/10/2023 12:02"
$ws.Range("P85").Value = 3.89
$ws.Range("Q85").Value = "04/11/2023 15:29"
$ws.Range("R85").Value = 3.8
$ws.Range("S85").Value = "22/10/2023 12:02"
$ws.Range("T85").Value = 3.52
$ws.Range("U85").Value = "04/11/2023 15:29"
$ws.Range("V85").Value = "https://www.betexplorer.com/football/germany/bundesliga/1-fc-koln-augsburg/CzIaCzO1/"

# Row 86
$ws.Range("A86").Value = 85
$ws.Range("B86").Value = "germany"
$ws.Range("C86").Value = "bundesliga"
$ws.Range("D86").Value = "2023-2024"
$ws.Range("E86").Value = 45234.64583333334
$ws.Range("F86").Value = "Freiburg"
$ws.Range("G86").Value = 3
$ws.Range("H86").Value = "B. Monchengladbach"
$ws.Range("I86").Value = 3
$ws.Range("J86").Value = 1.77
$ws.Range("K86").Value = "22/10/2023 12:02"
$ws.Range("L86").Value = 2.03
$ws.Range("M86").Value = "04/11/2023 15:19"
$ws.Range("N86").Value = 3.99
$ws.Range("O86").Value = "22/10/2023 12:02"
$ws.Range("P86").Value = 3.85
$ws.Range("Q86").Value = "04/11/2023 15:19"
$ws.Range("R86").Value = 4.55
$ws.Range("S86").Value = "22/10/2023 12:02"
$ws.Range("T86").Value = 3.64
$ws.Range("U86").Value = "04/11/2023 15:19"
$ws.Range("V86").Value = "https://www.betexplorer.com/football/germany/bundesliga/freiburg-b-monchengladbach/MoPnFxvq/"

# Row 87
$ws.Range("A87").Value = 86
$ws.Range("B87").Value = "germany"
$ws.Range("C87").Value = "bundesliga"
$ws.Range("D87").Value = "2023-2024"
$ws.Range("E87").Value = 45234.64583333334
$ws.Range("F87").Value = "Hoffenheim"
$ws.Range("G87").Value = 2
$ws.Range("H87").Value = "Bayer Leverkusen"
$ws.Range("I87").Value = 3
$ws.Range("J87").Value = 3.85
$ws.Range("K87").Value = "22/10/2023 12:02"
$ws.Range("L87").Value = 5.53
$ws.Range("M87").Value = "04/11/2023 15:29"
$ws.Range("N87").Value = 4.05
$ws.Range("O87").Value = "22/10/2023 12:02"
$ws.Range("P87").Value = 5.01
$ws.Range("Q87").Value = "04/11/2023 15:27"
$ws.Range("R87").Value = 1.81
$ws.Range("S87").Value = "22/10/2023 12:02"
$ws.Range("T87").Value = 1.53
$ws.Range("U87").Value = "04/11/2023 15:29"
$ws.Range("V87").Value = "https://www.betexplorer.com/football/germany/bundesliga/hoffenheim-bayer-leverkusen/jRG3Bfw8/"

# Row 88
$ws.Range("A88").Value = 87
$ws.Range("B88").Value = "germany"
$ws.Range("C88").Value = "bundesliga"
$ws.Range("D88").Value = "2023-2024"
$ws.Range("E88").Value = 45234.64583333334
$ws.Range("F88").Value = "Union Berlin"
$ws.Range("G88").Value = 0
$ws.Range("H88").Value = "Eintracht Frankfurt"
$ws.Range("I88").Value = 3
$ws.Range("J88").Value = 2.1
$ws.Range("K88").Value = "22/10/2023 12:02"
$ws.Range("L88").Value = 2.54
$ws.Range("M88").Value = "04/11/2023 15:28"
$ws.Range("N88").Value = 3.39
$ws.Range("O88").Value = "22/10/2023 12:02"
$ws.Range("P88").Value = 3.26
$ws.Range("Q88").Value = "04/11/2023 15:23"
$ws.Range("R88").Value = 3.8
$ws.Range("S88").Value = "22/10/2023 12:02"
$ws.Range("T88").Value = 3.06
$ws.Range("U88").Value = "04/11/2023 15:28"
$ws.Range("V88").Value = "https://www.betexplorer.com/football/germany/bundesliga/union-berlin-eintracht-frankfurt/ngYKJIVR/"

# Row 89
$ws.Range("A89").Value = 88
$ws.Range("B89").Value = "germany"
$ws.Range("C89").Value = "bundesliga"
$ws.Range("D89").Value = "2023-2024"
$ws.Range("E89").Value = 45234.77083333334
$ws.Range("F89").Value = "Dortmund"
$ws.Range("G89").Value = 0
$ws.Range("H89").Value = "Bayern Munich"
$ws.Range("I89").Value = 4
$ws.Range("J89").Value = 3.41
$ws.Range("K89").Value = "22/10/2023 12:02"
$ws.Range("L89").Value = 3.91
$ws.Range("M89").Value = "04/11/2023 18:29"
$ws.Range("N89").Value = 4.28
$ws.Range("O89").Value = "22/10/2023 12:02"
$ws.Range("P89").Value = 4.53
$ws.Range("Q89").Value = "04/11/2023 18:29"
$ws.Range("R89").Value = 1.88
$ws.Range("S89").Value = "22/10/2023 12:02"
$ws.Range("T89").Value = 1.81
$ws.Range("U89").Value = "04/11/2023 18:29"
$ws.Range("V89").Value = "https://www.betexplorer.com/football/germany/bundesliga/dortmund-bayern-munich/hWYGKbGL/"

